$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.045.98'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -2.96%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.025.03'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '547.38'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '135.68'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.019.18'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.15%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -4.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.09'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -5.10%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.450'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.32'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.514.74'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '62.122.87'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -2.99%  '
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.027.37'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.67'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '476.68'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.29'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.82%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.675'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -3.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.08'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '80.56'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.16'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.72'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.77'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -3.96%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.91'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '25.79'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.21%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.32'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '55.49'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.14%  '
$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.44'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.94'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '458.42'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -8.46%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.240.60'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0798'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0385'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.15'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.48'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -7.80%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '25.89'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +5.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.245'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -3.73%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.00'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.108'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '118.58'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -3.83%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₃0497'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -6.50%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +7.19%  '
